$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the product links in rows A3-A6 with new ones, and fill in the
# previously-empty A7 with a new link too.
#
# Final row -> URL mapping (matches the target workbook):
#   A3 -> Lenovo V15 monitor/laptop listing
#   A4 -> BenQ EX2710Q monitor
#   A5 -> Acer Swift laptops category page
#   A6 -> Asus TUF 34" curved gaming monitor
#   A7 -> Gigabyte graphics cards category page (new row content)

$ws.Range("A3").Value = "https://box.co.uk/82yu00hxuk-lenovo-v15-g4-amd-ryzen-5-16gb-256gb"
$ws.Range("A4").Value = "https://box.co.uk/9h-lk4la-tbe-benq-ex2710q-27in-2k-ultrahd-ips-gaming"
$ws.Range("A5").Value = "https://box.co.uk/acer-swift-laptops"
$ws.Range("A6").Value = "https://box.co.uk/vg34vqel1a-asus-tuf-34-uwqhd-led-100hz-curved-gaming"
$ws.Range("A7").Value = "https://box.co.uk/gigabyte-graphic-cards"

# Wire up the hyperlinks for each cell (in the same order the relationship
# ids were produced in the target file: A4, A5, A6, A3, A7).
$ws.Hyperlinks.Add($ws.Range("A4"), "https://box.co.uk/9h-lk4la-tbe-benq-ex2710q-27in-2k-ultrahd-ips-gaming")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://box.co.uk/acer-swift-laptops")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://box.co.uk/vg34vqel1a-asus-tuf-34-uwqhd-led-100hz-curved-gaming")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://box.co.uk/82yu00hxuk-lenovo-v15-g4-amd-ryzen-5-16gb-256gb")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://box.co.uk/gigabyte-graphic-cards")

# Adding the hyperlinks applies direct font formatting on top of the cells;
# restore the plain "Hyperlink" cell style (as used by the other linked
# cells/rows) so the cells keep looking consistent with A2.
$ws.Range("A3:A7").Style = "Hyperlink"

# Match the saved selection/active cell of the target file.
$ws.Range("A9").Select()
